$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure date-like strings in column H stay as text, not auto-converted dates
$ws.Range("H2:H19").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 'Issue'
$ws.Cells.Item(2, 3).Value = 'Deep'
$ws.Cells.Item(2, 4).Value = 'ShahDeep5113'
$ws.Cells.Item(2, 5).Value = 'open'
$ws.Cells.Item(2, 6).Value = 'ShahDeep5113'
$ws.Cells.Item(2, 7).Value = '1. Works on IPCL Report (Shell Breakage Report) In Database.
2. Create SPs for this report and updating templates.
3. GET Employee Master Table and used in the SBR.
4. testing Done for Shell Breakage '
$ws.Cells.Item(2, 8).Value = '2025-07-30'

# Row 3
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 'Issue'
$ws.Cells.Item(3, 3).Value = 'HONDA : TR : SR94291  For  YES/No and Date Validate'
$ws.Cells.Item(3, 4).Value = 'Vishal-Bhaliya'
$ws.Cells.Item(3, 5).Value = 'open'
$ws.Cells.Item(3, 6).Value = 'Vishal-Bhaliya'
$ws.Cells.Item(3, 7).Value = '<img width="903" height="587" alt="Image" src="https://github.com/user-attachments/assets/d1a6e77f-8d7a-4aad-838f-0aae1df7028c" />'
$ws.Cells.Item(3, 8).Value = '2025-07-30'

# Row 4
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 'Issue'
$ws.Cells.Item(4, 3).Value = 'HONDA : TR :- Batch/Employee Option'
$ws.Cells.Item(4, 4).Value = 'ShahDeep5113, Vishal-Bhaliya'
$ws.Cells.Item(4, 5).Value = 'open'
$ws.Cells.Item(4, 6).Value = 'Vishal-Bhaliya'
$ws.Cells.Item(4, 7).Value = '<img width="903" height="587" alt="Image" src="https://github.com/user-attachments/assets/863c314e-4191-4561-9c97-c0a6995e11a5" />'
$ws.Cells.Item(4, 8).Value = '2025-07-30'

# Row 5
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 'Issue'
$ws.Cells.Item(5, 3).Value = 'NP1141-INC-450913-add edit rights in SCR database for DTC team'
$ws.Cells.Item(5, 4).Value = 'Nandini-RI'
$ws.Cells.Item(5, 5).Value = 'closed'
$ws.Cells.Item(5, 6).Value = 'Nandini-RI'
$ws.Cells.Item(5, 7).Value = 'Assign all three people to "Editor" role.
Asked to user for verify.'
$ws.Cells.Item(5, 8).Value = '2025-07-30'

# Row 6
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 'Issue'
$ws.Cells.Item(6, 3).Value = 'Kontoor: Sarge - Re: SARGE ISSUE'
$ws.Cells.Item(6, 4).Value = 'Unassigned'
$ws.Cells.Item(6, 5).Value = 'open'
$ws.Cells.Item(6, 6).Value = 'Vrushali-gohel'
$ws.Cells.Item(6, 7).Value = '- Reset user account<M.Elghobary@lotusgarments.com> password in sarge production'
$ws.Cells.Item(6, 8).Value = '2025-07-30'

# Row 7
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 'Issue'
$ws.Cells.Item(7, 3).Value = 'Kontoor: SCR0217 - InfoCenter CoC: Edit function amendment []'
$ws.Cells.Item(7, 4).Value = 'Kinjal-Makwana, Ravi-Morichauhan'
$ws.Cells.Item(7, 5).Value = 'open'
$ws.Cells.Item(7, 6).Value = 'Vrushali-gohel'
$ws.Cells.Item(7, 7).Value = '- Apply logic to edit maker data and test.'
$ws.Cells.Item(7, 8).Value = '2025-07-30'

# Row 8
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 'Issue'
$ws.Cells.Item(8, 3).Value = 'Vrushali'
$ws.Cells.Item(8, 4).Value = 'Vrushali-gohel'
$ws.Cells.Item(8, 5).Value = 'open'
$ws.Cells.Item(8, 6).Value = 'Vrushali-gohel'
$ws.Cells.Item(8, 7).Value = '- https://github.com/RI-BVN/RamansheeRepo/issues/230 -- Checked log
- https://github.com/RI-BVN/RamansheeRepo/issues/633 -- reset account password
- https://github.com/RI-BVN/RamansheeRepo/issues/450 '
$ws.Cells.Item(8, 8).Value = '2025-07-30'

# Row 9
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 'Issue'
$ws.Cells.Item(9, 3).Value = 'Miloni'
$ws.Cells.Item(9, 4).Value = 'Miloni-Mehta04'
$ws.Cells.Item(9, 5).Value = 'open'
$ws.Cells.Item(9, 6).Value = 'Miloni-Mehta04'
$ws.Cells.Item(9, 7).Value = 'Working on UI for VMI CRS Account Function and some modification in it'
$ws.Cells.Item(9, 8).Value = '2025-07-30'

# Row 10
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 'Issue'
$ws.Cells.Item(10, 3).Value = 'Ravi'
$ws.Cells.Item(10, 4).Value = 'Ravi-Morichauhan'
$ws.Cells.Item(10, 5).Value = 'open'
$ws.Cells.Item(10, 6).Value = 'Ravi-Morichauhan'
$ws.Cells.Item(10, 7).Value = '1. RoofTopLight Company Stock List – Ensure proper filtering, pagination, and design.
2. RoofTopLight ConfirmOrder - proper filtering, pagination, and design.
3. RoofTopLight Challan List - proper fil'
$ws.Cells.Item(10, 8).Value = '2025-07-30'

# Row 11
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 'Issue'
$ws.Cells.Item(11, 3).Value = 'Utsav'
$ws.Cells.Item(11, 4).Value = 'Utsav-Patel9'
$ws.Cells.Item(11, 5).Value = 'open'
$ws.Cells.Item(11, 6).Value = 'Utsav-Patel9'
$ws.Cells.Item(11, 7).Value = 'Powerapp :
1. Daily Stock Adjustment List in second screen in add field
2. Daily Stock Adjustment List in design complete
3. Daily Stock Adjustment List in click on edit and all screen value bind scre'
$ws.Cells.Item(11, 8).Value = '2025-07-30'

# Row 12
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 'Issue'
$ws.Cells.Item(12, 3).Value = 'Rakesh'
$ws.Cells.Item(12, 4).Value = 'Rakesh-Morichauhan'
$ws.Cells.Item(12, 5).Value = 'open'
$ws.Cells.Item(12, 6).Value = 'Rakesh-Morichauhan'
$ws.Cells.Item(12, 7).Value = ' 
1. Add/Update (Alert Massage in Insert successfully completed Or Update Successfully Completed) 
2. Insert Data Receipt Upload (Optional)  jpg, png, etc.
3. Department Head List page and Edit models'
$ws.Cells.Item(12, 8).Value = '2025-07-30'

# Row 13
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 'Issue'
$ws.Cells.Item(13, 3).Value = 'HONDA : SR101180 :Travel billing module'
$ws.Cells.Item(13, 4).Value = 'ShahDeep5113, Vishal-Bhaliya, Urmi-Parmar'
$ws.Cells.Item(13, 5).Value = 'open'
$ws.Cells.Item(13, 6).Value = 'Vishal-Bhaliya'
$ws.Cells.Item(13, 7).Value = '<img width="800" height="458" alt="Image" src="https://github.com/user-attachments/assets/16cf95b2-3d43-4377-8cf6-0d67a4c3fb0f" />'
$ws.Cells.Item(13, 8).Value = '2025-07-30'

# Row 14
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = 'Issue'
$ws.Cells.Item(14, 3).Value = 'Kinjal'
$ws.Cells.Item(14, 4).Value = 'Kinjal-Makwana'
$ws.Cells.Item(14, 5).Value = 'open'
$ws.Cells.Item(14, 6).Value = 'Kinjal-Makwana'
$ws.Cells.Item(14, 7).Value = '**In 1st Half**
- Ticket, Outstanding file update
- #617 - Replied user
- #590 - Replied user
- #618 - Replied user
- Weekly meeting from 1 to 2
**In 2nd Half**
- #630 - Work and replied user
- https'
$ws.Cells.Item(14, 8).Value = '2025-07-30'

# Row 15
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 'Issue'
$ws.Cells.Item(15, 3).Value = 'Riddhi - TIR Prod SAP Auto Export Status'
$ws.Cells.Item(15, 4).Value = 'Vrushali-gohel, RiddhiBaraiya'
$ws.Cells.Item(15, 5).Value = 'open'
$ws.Cells.Item(15, 6).Value = 'Vrushali-gohel'
$ws.Cells.Item(15, 7).Value = '<img width="1218" height="560" alt="Image" src="https://github.com/user-attachments/assets/daedf98b-d8ac-4c6e-984f-f702cd561b79" />'
$ws.Cells.Item(15, 8).Value = '2025-07-31'

# Row 16
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = 'Issue'
$ws.Cells.Item(16, 3).Value = 'HONDA : SR102099 -Emp.Code-12782-Rellocation click forgated'
$ws.Cells.Item(16, 4).Value = 'Vishal-Bhaliya'
$ws.Cells.Item(16, 5).Value = 'open'
$ws.Cells.Item(16, 6).Value = 'Vishal-Bhaliya'
$ws.Cells.Item(16, 7).Value = '<img width="900" height="757" alt="Image" src="https://github.com/user-attachments/assets/2afa467f-943d-45a1-97ad-87fc1694243b" />'
$ws.Cells.Item(16, 8).Value = '2025-07-31'

# Row 17
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = 'Issue'
$ws.Cells.Item(17, 3).Value = 'HONDA : PR dashboard access - SR99647'
$ws.Cells.Item(17, 4).Value = 'Urmi-Parmar'
$ws.Cells.Item(17, 5).Value = 'open'
$ws.Cells.Item(17, 6).Value = 'Vishal-Bhaliya'
$ws.Cells.Item(17, 7).Value = 'change log shared with  <img width="800" height="458" alt="Image" src="https://github.com/user-attachments/assets/7a420465-48f4-4980-838e-4da8ed0cd776" />'
$ws.Cells.Item(17, 8).Value = '2025-07-31'

# Row 18
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = 'Issue'
$ws.Cells.Item(18, 3).Value = 'HONDA : SR100081 : Roaster Extra Seat add For Div  based on OP seat margin'
$ws.Cells.Item(18, 4).Value = 'ShahDeep5113, Vishal-Bhaliya'
$ws.Cells.Item(18, 5).Value = 'open'
$ws.Cells.Item(18, 6).Value = 'Vishal-Bhaliya'
$ws.Cells.Item(18, 7).Value = 'Need To correct Some Functionality after discussion with satyaveer san aS ON 30-07-205'
$ws.Cells.Item(18, 8).Value = '2025-07-31'

# Row 19
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = 'Issue'
$ws.Cells.Item(19, 3).Value = 'HONDA : SR104589 - SHE Documents on Employee portal'
$ws.Cells.Item(19, 4).Value = 'Vishal-Bhaliya'
$ws.Cells.Item(19, 5).Value = 'closed'
$ws.Cells.Item(19, 6).Value = 'Vishal-Bhaliya'
$ws.Cells.Item(19, 7).Value = 'Go live on 30-07-2025 3:30 PM '
$ws.Cells.Item(19, 8).Value = '2025-07-31'
